# Rename "shortTrainingPost1998" to "trainingMeasures" in cell A14,
# and update the saved view/selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the program identifier value in A14.
$ws.Range("A14").Value = "trainingMeasures"

# Update the sheet's scroll position / selection to match the saved view.
$ws.Activate()
$ws.Range("A15").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
